$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume table with the latest scraped values.
# Row 42 and 43 swap places in ranking (dogwifhat now above Bittensor), so
# all four columns for those two rows are updated accordingly.
# Note: price values in column D are text (not numbers) in this sheet - a
# leading apostrophe forces Excel to store them as text rather than
# auto-converting to a numeric value.
$ws.Range("D2").Value = "68.832.37"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.923.02"
$ws.Range("E3").Value = "  +4.53%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'605.04"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'165.68"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "3.920.90"
$ws.Range("E7").Value = "  +4.56%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("D11").Value = "'6.41"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "'37.31"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "'0.0000246"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "4.582.01"
$ws.Range("E15").Value = "  +4.66%  "
$ws.Range("D16").Value = "3.939.19"
$ws.Range("E16").Value = "  +5.21%  "
$ws.Range("D17").Value = "68.982.94"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").Value = "'11.16"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").Value = "'487.79"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "'0.723"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  +12.23%  "
$ws.Range("D25").Value = "'84.47"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D27").Value = "'12.11"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").Value = "'10.16"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'2.95"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "4.077.03"
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'7.86"
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("D34").Value = "'32.27"
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("D35").Value = "3.874.13"
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").Value = "'5.92"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.01"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'437.92"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'48.47"
$ws.Range("D46").Value = "'8.51"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "2.851.77"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").Value = "'26.27"
$ws.Range("E49").Value = "  +11.28%  "
$ws.Range("D50").Value = "'142.08"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +1.37%  "
